# Apply updated "dSF" (column F) values pulled from the latest data source.
# Maps worksheet row number -> new value for column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    9  = 0
    11 = -4
    13 = -1
    20 = -2
    21 = -9
    24 = -3
    27 = 3
    34 = 3
    36 = -3
    43 = -6
    44 = -3
    46 = -2
    47 = 5
    49 = 4
    50 = -2
    51 = -8
    52 = -9
    53 = -6
    54 = -2
    56 = -7
    58 = -3
    59 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
